$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '87.766.23'
$ws.Range('E2').Value = '  -2.39%  '
$ws.Range('D3').Value = '3.049.84'
$ws.Range('E3').Value = '  -5.29%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '207.59'
$ws.Range('E5').Value = '  -5.12%  '
$ws.Range('D6').Value = '613.67'
$ws.Range('E6').Value = '  -2.49%  '
$ws.Range('D7').Value = '0.362'
$ws.Range('E7').Value = '  -7.72%  '
$ws.Range('D8').Value = '0.803'
$ws.Range('E8').Value = '  +15.22%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').Value = '3.043.96'
$ws.Range('E10').Value = '  -5.37%  '
$ws.Range('D11').Value = '0.589'
$ws.Range('E11').Value = '  +2.82%  '
$ws.Range('E12').Value = '  -1.05%  '
$ws.Range('D13').Value = '0.0000232'
$ws.Range('E13').Value = '  -10.64%  '
$ws.Range('E14').Value = '  -4.00%  '
$ws.Range('D15').Value = '87.382.93'
$ws.Range('E15').Value = '  -2.46%  '
$ws.Range('D16').Value = '3.593.78'
$ws.Range('E16').Value = '  -5.82%  '
$ws.Range('D17').Value = '31.22'
$ws.Range('E17').Value = '  -6.85%  '
$ws.Range('D18').Value = '3.012.77'
$ws.Range('E18').Value = '  -5.79%  '
$ws.Range('D19').Value = '3.15'
$ws.Range('E19').Value = '  -9.44%  '
$ws.Range('B20').Value = 'PEPE'
$ws.Range('C20').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D20').Value = '0.0000194'
$ws.Range('E20').Value = '  -17.40%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').Value = '12.97'
$ws.Range('E21').Value = '  -4.45%  '
$ws.Range('D22').Value = '415.73'
$ws.Range('E22').Value = '  -5.53%  '
$ws.Range('D23').Value = '7.99'
$ws.Range('E23').Value = '  -7.59%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '4.80'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -6.00%  '
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('D26').Value = '11.49'
$ws.Range('E26').Value = '  -3.85%  '
$ws.Range('D27').Value = '80.27'
$ws.Range('E27').Value = '  -1.51%  '
$ws.Range('D28').Value = '3.225.78'
$ws.Range('E28').Value = '  -4.85%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +0.40%  '
$ws.Range('E30').Value = '  +8.55%  '
$ws.Range('D31').Value = '0.159'
$ws.Range('E31').Value = '  -0.86%  '
$ws.Range('D32').Value = '7.95'
$ws.Range('E32').Value = '  -7.41%  '
$ws.Range('D33').Value = '497.66'
$ws.Range('E33').Value = '  -8.63%  '
$ws.Range('D34').Value = '3.46'
$ws.Range('E34').Value = '  -16.52%  '
$ws.Range('D35').Value = '6.48'
$ws.Range('E35').Value = '  -7.75%  '
$ws.Range('E36').Value = '  -7.65%  '
$ws.Range('E37').Value = '  -7.32%  '
$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '21.90'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -2.39%  '
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').Value = '22.22'
$ws.Range('E39').Value = '  -0.71%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = '0.129'
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').Value = '0.354'
$ws.Range('E43').Value = '  -5.80%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').Value = '146.56'
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('E45').Value = '  +6.52%  '
$ws.Range('D46').Value = '1.77'
$ws.Range('E46').Value = '  -8.26%  '
$ws.Range('D47').Value = '43.18'
$ws.Range('E47').Value = '  -1.40%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.0660'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +8.12%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '153.63'
$ws.Range('E49').Value = '  -11.60%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '0.691'
$ws.Range('E50').Value = '  -7.51%  '
$ws.Range('D51').Value = '1.16'
$ws.Range('E51').Value = '  -8.29%  '
